$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 4 = "Table Level Check Constraint"
# Column 2 ("Quantity"): change highlighted "4" to "3"
$qtyCell = $t.Cell(4, 2)
$qtyCell.Range.Text = "3"

# Column 3 ("Comments"): add explanatory comment to the previously empty paragraph
$commentCell = $t.Cell(4, 3)
$commentCell.Range.Text = "Added 3 table level check constraints"
